# Add 2022-Q4 data: insert a new "2022-Q4" worksheet (fund holdings detail,
# copied/styled after the existing "2022-Q3" sheet) and add the
# corresponding summary row to the "总计" (Total) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q4" and shift the
#    existing quarters down by one row.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Re-apply formatting to the freshly inserted row: column A keeps the
# bordered/bold "index" style (same as the rows below), B:D stay
# unformatted like the other data rows.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.12

# Bump the running index of the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# ---------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right before "2022-Q3" (so the
#    tab order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2, ...).
# ---------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($sheetQ3)
$newSheet.Name = "2022-Q4"

# Re-fetch fresh references after the rename.
$src = $wb.Worksheets.Item("2022-Q3")
$dst = $wb.Worksheets.Item("2022-Q4")

# Clone the header/data-row formatting (bold+bordered header, bordered
# index column) from the existing "2022-Q3" sheet.
$src.Range("A1:H2").Copy()
$dst.Range("A1").PasteSpecial(-4122)
$dst.Range("A2:H2").Copy()
$dst.Range("A3:H6").PasteSpecial(-4122)

# Header row
$dst.Range("B1").Value = "基金代码"
$dst.Range("C1").Value = "基金名称"
$dst.Range("D1").Value = "基金规模"
$dst.Range("E1").Value = "股票总仓位"
$dst.Range("F1").Value = "仓位占比"
$dst.Range("G1").Value = "持有市值(亿元)"
$dst.Range("H1").Value = "仓位排名"

# Fund holdings data. B (code) and D:G (scale/position/ratio/value) are
# stored as plain text in this workbook (leading apostrophe keeps them
# text instead of being auto-converted to numbers); H (rank) is numeric,
# same as the other quarter sheets.
$rows = @(
  @{A=0; B="010495"; C="创金合信创新驱动股票A"; D="0.93"; E="90.74"; F="5.03"; G="0.0468"; H=5},
  @{A=1; B="001972"; C="前海开源沪港深智慧生活优选灵活配置混合"; D="0.52"; E="91.71"; F="6.39"; G="0.0332"; H=8},
  @{A=2; B="010496"; C="创金合信创新驱动股票C"; D="0.32"; E="90.74"; F="5.03"; G="0.0161"; H=5},
  @{A=3; B="009569"; C="浙商智多宝稳健一年持有期混合C"; D="1.02"; E="26.91"; F="1.42"; G="0.0145"; H=2},
  @{A=4; B="009568"; C="浙商智多宝稳健一年持有期混合A"; D="0.98"; E="26.91"; F="1.42"; G="0.0139"; H=2}
)

$r = 2
foreach ($row in $rows) {
  $dst.Cells.Item($r, 1).Value = $row.A
  $dst.Cells.Item($r, 2).Value = "'" + $row.B
  $dst.Cells.Item($r, 3).Value = $row.C
  $dst.Cells.Item($r, 4).Value = "'" + $row.D
  $dst.Cells.Item($r, 5).Value = "'" + $row.E
  $dst.Cells.Item($r, 6).Value = "'" + $row.F
  $dst.Cells.Item($r, 7).Value = "'" + $row.G
  $dst.Cells.Item($r, 8).Value = $row.H
  $r = $r + 1
}

# Restore the original active sheet ("总计") so the workbook's view state
# isn't disturbed by the sheet-insertion/activation side effects above.
$total.Activate()
